# Rearranges (permutes) the species-observation rows 15-28 and 30-38 on the
# active sheet so that each target row ends up holding the data that
# originally belonged to a different row within that same block (row 29 is
# left untouched). Only the columns whose values actually differ from one
# record to another are touched (A, B, D, E, F, G, H, M, P, Q, R); the many
# columns that are identical for every record in the block (dates, booleans,
# location/county text, empty placeholder cells, etc.) are left completely
# alone so no accidental type/format conversions are introduced.
#
# $mapping[newRow] = oldRow  -> newRow must end up with the content that
# originally lived in oldRow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    15 = 21
    16 = 20
    17 = 32
    18 = 16
    19 = 35
    20 = 36
    21 = 22
    22 = 33
    23 = 31
    24 = 26
    25 = 30
    26 = 23
    27 = 37
    28 = 18
    30 = 34
    31 = 25
    32 = 19
    33 = 17
    34 = 27
    35 = 38
    36 = 24
    37 = 28
    38 = 15
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "P", "Q", "R")

# 1. Snapshot the values of the columns that vary, for every source row,
#    *before* any writes happen, so overlapping cycles in the permutation
#    don't clobber data that is still needed later on.
$snapshots = @{}
foreach ($oldRow in $mapping.Values) {
    if (-not $snapshots.ContainsKey($oldRow)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$oldRow").Value2
        }
        $snapshots[$oldRow] = $rowData
    }
}

# 2. Write each snapshot into its new row, column by column.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $rowData = $snapshots[$oldRow]
    foreach ($col in $cols) {
        $val = $rowData[$col]
        if ($null -eq $val) { $val = "" }
        $ws.Range("$col$newRow").Value2 = $val
    }
}
